# Refresh the cryptos price/volume snapshot (GitHub Actions data pull).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.784.74"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "'1.638.35"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'215.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("E6").Value = "  -0.55%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").Value = "'19.64"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "'1.864.78"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "
$ws.Range("D14").Value = "'1.641.54"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.50%  "
$ws.Range("D15").Value = "'0.563"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "'63.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.27%  "
$ws.Range("D18").Value = "'25.831.72"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "'192.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "'9.95"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("D24").Value = "'1.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.38%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'141.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("E27").Value = "  +1.55%  "
$ws.Range("E28").Value = "  +1.57%  "
$ws.Range("E29").Value = "  -0.07%  "
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").Value = "'0.0493"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("E32").Value = "  +0.82%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("E36").Value = "  +0.48%  "
$ws.Range("D37").Value = "'1.135.15"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.26%  "
$ws.Range("E38").Value = "  -1.62%  "
$ws.Range("E39").Value = "  -0.93%  "
$ws.Range("E40").Value = "  +0.06%  "
$ws.Range("E41").Value = "  +0.19%  "
$ws.Range("E42").Value = "  +0.72%  "
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("E44").Value = "  +0.62%  "
$ws.Range("D45").Value = "'1.774.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.16%  "
$ws.Range("E46").Value = "  +2.67%  "
$ws.Range("D47").Value = "'55.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.32%  "
$ws.Range("E48").Value = "  -1.17%  "
$ws.Range("D50").Value = "'1.42"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("E51").Value = "  -2.55%  "
